# Append the new daily allocation row (2025-11-18) to the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 77 -> 78).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Column A holds dates formatted as plain text (e.g. "09/03/2025"), not real
# date serials. A leading apostrophe forces Excel to store the literal text
# instead of auto-converting it to a date value; we then reset the cell's
# style back to Normal so no stray quote-prefixed number format lingers on
# the cell (keeping it identical in shape to the other data rows).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = "'11/18/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 0.2054993753108948
$ws.Cells.Item($newRow, 3).Value = 0.7945006246891052
